# Updated cryptos list on Sat Jun 29 14:59:27 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures with the
# latest crypto-ranking scrape. Two coins (NEARProtocol and
# RenzoRestakedETH) also swapped rank positions (rows 37/38), so their
# name/link/price/volume cells are rewritten in place.
#
# Column D values are plain digits/dots (e.g. "571.29", "1.00"), which
# Excel would otherwise auto-convert to numbers (dropping the authored
# text formatting, e.g. "1.00" -> 1). Forcing the cell to Text format
# ("@") before the write keeps them as text, matching how this sheet's
# price column was already authored. Column E values already contain
# padding spaces and a "%" sign, and column B/C are clearly non-numeric,
# so they round-trip as text without needing that treatment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.983.28"
$ws.Range("E2").Value = "  +0.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.386.44"
$ws.Range("E3").Value = "  -0.64%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.29"
$ws.Range("E5").Value = "  -0.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.81"
$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  -0.66%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.61"
$ws.Range("E9").Value = "  +1.01%  "

$ws.Range("E10").Value = "  -1.42%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.390"
$ws.Range("E11").Value = "  -0.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.966.10"
$ws.Range("E12").Value = "  -0.61%  "

$ws.Range("E13").Value = "  +2.21%  "

$ws.Range("E14").Value = "  -2.04%  "

$ws.Range("E15").Value = "  -0.83%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.383.34"
$ws.Range("E16").Value = "  -0.94%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.080.49"
$ws.Range("E17").Value = "  +0.19%  "

$ws.Range("E18").Value = "  -3.60%  "

$ws.Range("E19").Value = "  -4.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.91"
$ws.Range("E20").Value = "  -4.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "381.85"
$ws.Range("E21").Value = "  -3.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.91"
$ws.Range("E22").Value = "  +2.63%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.552"
$ws.Range("E23").Value = "  -2.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.80%  "

$ws.Range("E25").Value = "  -4.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.523.05"
$ws.Range("E26").Value = "  -1.16%  "

$ws.Range("E27").Value = "  +1.24%  "

$ws.Range("E28").Value = "  +0.44%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.30"
$ws.Range("E29").Value = "  -2.33%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.16"
$ws.Range("E30").Value = "  -0.51%  "

$ws.Range("E31").Value = "  -2.45%  "

$ws.Range("E32").Value = "  -2.46%  "

$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.35"
$ws.Range("E34").Value = "  -2.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.95"
$ws.Range("E35").Value = "  -0.54%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "166.22"
$ws.Range("E36").Value = "  -0.63%  "

# Rows 37/38 swapped ranking position: NEARProtocol <-> RenzoRestakedETH,
# each carrying its own refreshed price/volume figures.
$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.418.48"
$ws.Range("E37").Value = "  -0.48%  "

$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.01"
$ws.Range("E38").Value = "  -2.46%  "

$ws.Range("E39").Value = "  -4.95%  "

$ws.Range("E40").Value = "  -2.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.06"
$ws.Range("E41").Value = "  -0.45%  "

$ws.Range("E42").Value = "  -0.11%  "

$ws.Range("E43").Value = "  -1.84%  "

$ws.Range("E44").Value = "  -2.48%  "

$ws.Range("E45").Value = "  -2.94%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.13"
$ws.Range("E46").Value = "  -0.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.457.64"
$ws.Range("E47").Value = "  -5.43%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.06"
$ws.Range("E48").Value = "  +0.49%  "

$ws.Range("E49").Value = "  -2.97%  "

$ws.Range("E50").Value = "  +2.09%  "

$ws.Range("E51").Value = "  +7.36%  "
